$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: updated per-epoch accuracy values (re-run after freezing
#     token embeddings + decoder layer 1) ---
$accByRow = @{
    4=0.640625; 5=0.625; 6=0.53125; 7=0.515625; 8=0.4375; 9=0.515625;
    10=0.53125; 11=0.546875; 12=0.5625; 13=0.53125; 14=0.484375; 15=0.46875;
    16=0.53125; 17=0.5; 18=0.421875; 19=0.453125; 20=0.4375; 21=0.421875;
    22=0.421875; 23=0.40625;
    31=0.40625; 32=0.40625; 33=0.40625; 34=0.40625;
    43=0.390625;
    54=0.390625; 55=0.390625; 56=0.390625; 57=0.390625; 58=0.390625;
    59=0.390625; 60=0.390625; 61=0.390625; 62=0.390625; 63=0.390625;
    64=0.390625; 65=0.390625; 66=0.390625; 67=0.390625; 68=0.390625;
    69=0.390625; 70=0.390625; 71=0.390625; 72=0.390625; 73=0.390625;
    74=0.390625; 75=0.390625; 76=0.390625;
    77=0.40625; 78=0.40625; 79=0.40625; 80=0.40625; 81=0.40625; 82=0.40625;
    83=0.40625; 84=0.40625; 85=0.40625; 86=0.40625; 87=0.40625; 88=0.40625;
    89=0.40625; 90=0.40625; 91=0.40625; 92=0.40625; 93=0.40625; 94=0.40625;
    95=0.40625; 96=0.40625; 97=0.40625; 98=0.40625; 99=0.40625; 100=0.40625;
    101=0.40625; 102=0.40625;
    103=0.421875; 104=0.28125; 105=0.4375; 106=0.328125; 107=0.296875;
    108=0.359375; 109=0.484375; 110=0.4375; 111=0.4375; 112=0.5625;
    113=0.453125; 114=0.28125; 115=0.375; 116=0.4375; 117=0.296875
}

foreach ($row in $accByRow.Keys) {
    $ws.Cells.Item($row, 2).Value = $accByRow[$row]
}

# --- Column A (rows 102-118): the repr()'d Python object id baked into the
#     dump moved since the notebook kernel restarted for this run ---
$newRepr = "<__main__.DisplayOutputs object at 0x7f36601d0f40>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newRepr
}

# --- Restore the active-cell/selection left by the notebook export ---
$ws.Range("A1:XFD1048576").Select()
$ws.Range("O13").Activate()
